$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: Last Updated timestamp (A2) ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 12:20 PM"

# --- Update "Industry Analysis" sheet: column F ("1 Year") values ---
$ws = $wb.Worksheets.Item("Industry Analysis")

$data = @{
    2 = 21.0016
    3 = -16.2396
    4 = 27.1317
    5 = -50.6494
    6 = 53.2813
    7 = -8.106199999999999
    8 = -9.552099999999999
    9 = 36.3756
    10 = -6.1314
    11 = 31.9081
    12 = -18.4955
    13 = 14.0155
    14 = -36.0718
    15 = -0.1622
    16 = 0.1459
    17 = -22.0012
    18 = 1.0561
    19 = -27.708
    20 = 47.7309
    21 = 12.0959
    22 = 95.1491
    23 = -50.2657
    24 = -13.3427
    25 = -9.9316
    26 = 5.8244
    27 = -32.7692
    28 = -24.8224
    29 = -18.4191
    30 = 25.8569
    31 = 58.4712
    32 = -3.3862
    33 = -6.3282
    34 = 27.7203
    35 = 4.4873
    36 = -4.9458
    37 = 3.6074
    38 = -23.3973
    39 = 8.7355
    40 = -5.8541
    41 = -8.3934
    42 = 20.3818
    43 = 14.3164
    44 = -12.6846
    45 = 28.4075
    46 = -1.1135
    47 = -37.1997
    48 = -29.8569
    49 = -27.5511
    50 = -49.7478
    51 = -51.8002
    52 = -38.5254
    53 = -12.4886
    54 = -5.0725
    55 = -17.7445
    56 = -26.636
    57 = -29.3361
    58 = -11.9574
    59 = -24.5687
    60 = -12.3
    61 = -10.9446
    62 = -17.1229
    63 = -9.5038
    64 = 54.2749
    65 = -43.4736
    66 = 13.2687
    67 = 12.7149
    68 = 24.8057
    69 = -17.0328
    70 = -6.8927
    71 = 13.6034
    72 = 3.9995
    73 = -16.226
    74 = -16.2448
    75 = 28.6924
    76 = 48.9752
}

foreach ($row in $data.Keys) {
    $ws.Cells.Item($row, 6).Value = $data[$row]
}
